# Add a new single-column table "xlsx_single_col_table" over L21:L25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header + data for the new table column.
$ws.Range("L21").Value = "colA"
$ws.Range("L22").Value = "a"
$ws.Range("L23").Value = "b"
$ws.Range("L24").Value = "c"
$ws.Range("L25").Value = "d"

# Create the table (ListObject) from the range, with a header row.
$null = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("L21:L25"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)

# Rename it (look it up by its auto-generated default name to avoid any
# stale handle returned from Add()).
$newTable = $ws.ListObjects.Item("Table4")
$newTable.Name = "xlsx_single_col_table"

# Match the new selection recorded in the workbook.
$null = $ws.Range("L22").Select()
